# Add an "Sr No" column (D) to the left of the existing table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "Sr No"
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 2
$ws.Range("D10").Value = 3
$ws.Range("D11").Value = 4
$ws.Range("D12").Value = 5

# Match the author's final selection/view state.
$ws.Range("D13").Select() | Out-Null
